$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Samsung UN65NU6900FXZA 65 inch 2160p LED Smart TV"

$ws.Range("D2").Formula = "=""579.99$"""
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)

$ws.Columns.Item(2).ColumnWidth = 47.14

$ws.Range("F8").Select()
